$wb = $excel.ActiveWorkbook

$wsRunManager = $wb.Worksheets.Item(1)
$wsData = $wb.Worksheets.Item(2)

# Update "execute" values on RUNMANAGER sheet (column C): no -> yes
$wsRunManager.Range("C2").Value = "yes"
$wsRunManager.Range("C3").Value = "yes"

# Update "execute" values on DATA sheet (column B)
$wsData.Range("B2").Value = "no"
$wsData.Range("B3").Value = "no"
$wsData.Range("B4").Value = "yes"
$wsData.Range("B6").Value = "yes"

# Update sheet selections / active view state:
# DATA sheet: selection moves from A1:H1 to B4, and it is no longer the tab shown when opening
$wsData.Activate()
$wsData.Range("B4").Select()

# RUNMANAGER sheet becomes the active/selected tab (its own selection C4 stays unchanged)
$wsRunManager.Activate()
